$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update grading values in column G
$ws.Range("G4").Value = 0
$ws.Range("G8").Value = 1
$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 1
$ws.Range("G11").Value = 1
$ws.Range("G12").Value = 6
$ws.Range("G15").Value = 3
$ws.Range("G16").Value = 0
$ws.Range("G18").Value = 0

# Clear the "Extra" section values entirely (cells removed from the sheet)
$ws.Range("G24").ClearContents()
$ws.Range("G25").ClearContents()
$ws.Range("G26").ClearContents()
$ws.Range("G27").ClearContents()
$ws.Range("G28").ClearContents()

# Update the active selection to match the saved view state
$ws.Range("G19").Select()

$wb.Save()
